$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text-typed value into a cell without leaving a stray number-format style behind.
# Excel treats digit-only / numeric-looking strings as numbers unless the cell is explicitly
# formatted as Text first; we then restore the style so no visible formatting change remains.
function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# --- Row 5 existing values -> text-formatted numeric strings ---
Set-TextCell 'F5' '130.00'
Set-TextCell 'G5' '2'
Set-TextCell 'H5' '2'
Set-TextCell 'I5' '130.00'
Set-TextCell 'J5' '65.00'

# --- Row 6 existing values -> text-formatted numeric strings ---
Set-TextCell 'F6' '40.00'
Set-TextCell 'G6' '1'
Set-TextCell 'H6' '1'
Set-TextCell 'I6' '40.00'
Set-TextCell 'J6' '40.00'

# --- New rows 7-15 ---
# Row 7
$ws.Range('A7').Value = 3
$ws.Range('B7').Value = 'wurth'
$ws.Range('C7').Value = '25-06-06'
$ws.Range('D7').Value = 'Sealing ring, vehicle RG-SEAL-DIN7603-HONDA-ALU-14X22X2'
Set-TextCell 'E7' '046411422'
Set-TextCell 'F7' '48.16'
Set-TextCell 'G7' '25'
Set-TextCell 'H7' '25'
Set-TextCell 'I7' '12.04'
$ws.Range('J7').Value = 0.48
Set-TextCell 'K7' '850 %'
$ws.Range('L7').Value = 4.56
$ws.Range('M7').Value = 0.46
$ws.Range('N7').Value = 5.02

# Row 8
$ws.Range('A8').Value = 4
$ws.Range('B8').Value = 'wurth'
$ws.Range('C8').Value = '25-06-06'
$ws.Range('D8').Value = 'Copper sealing ring RG-SEAL-DIN7603-CU-C-ASBESTFREE-20X26'
Set-TextCell 'E8' '046202026'
Set-TextCell 'F8' '324.82'
Set-TextCell 'G8' '20'
Set-TextCell 'H8' '20'
Set-TextCell 'I8' '64.96'
$ws.Range('J8').Value = 3.25
Set-TextCell 'K8' '850 %'
$ws.Range('L8').Value = 30.87
$ws.Range('M8').Value = 3.09
$ws.Range('N8').Value = 33.96

# Row 9
$ws.Range('A9').Value = 5
$ws.Range('B9').Value = 'wurth'
$ws.Range('C9').Value = '25-06-06'
$ws.Range('D9').Value = 'Radiator anti-freeze RADANTIFREZ-OAT-RED-20LTR'
Set-TextCell 'E9' '0892352020'
Set-TextCell 'F9' '108.00'
Set-TextCell 'G9' '1'
Set-TextCell 'H9' '1'
Set-TextCell 'I9' '108.01'
$ws.Range('J9').Value = 108.01
Set-TextCell 'K9' '850 %'
$ws.Range('L9').Value = 1026.1
$ws.Range('M9').Value = 102.61
$ws.Range('N9').Value = 1128.71

# Row 10
$ws.Range('A10').Value = 6
$ws.Range('B10').Value = 'wurth'
$ws.Range('C10').Value = '25-06-06'
$ws.Range('D10').Value = 'Halogen bulb, vehicle BULB-(H9005-HB3)-P20D-12V-60W'
Set-TextCell 'E10' '07201103'
Set-TextCell 'F10' '15.45'
Set-TextCell 'G10' '5'
Set-TextCell 'H10' '5'
Set-TextCell 'I10' '77.25'
$ws.Range('J10').Value = 15.45
Set-TextCell 'K10' '850 %'
$ws.Range('L10').Value = 146.77
$ws.Range('M10').Value = 14.68
$ws.Range('N10').Value = 161.45

# Row 11
$ws.Range('A11').Value = 7
$ws.Range('B11').Value = 'wurth'
$ws.Range('C11').Value = '25-06-06'
$ws.Range('D11').Value = 'Indicator/brake light bulb, vehicle BULB-INDCTR/BRK-P21/5W-BAY15D-12V-21/5W'
Set-TextCell 'E11' '07201341'
Set-TextCell 'F11' '2.87'
Set-TextCell 'G11' '20'
Set-TextCell 'H11' '20'
Set-TextCell 'I11' '57.40'
$ws.Range('J11').Value = 2.87
Set-TextCell 'K11' '850 %'
$ws.Range('L11').Value = 27.26
$ws.Range('M11').Value = 2.73
$ws.Range('N11').Value = 29.99

# Row 12
$ws.Range('A12').Value = 8
$ws.Range('B12').Value = 'wurth'
$ws.Range('C12').Value = '25-06-06'
$ws.Range('D12').Value = 'Silicone spray, universal CURTAIN-TRACK-LUBE-300G'
Set-TextCell 'E12' '0893125080'
Set-TextCell 'F12' '12.95'
Set-TextCell 'G12' '1'
Set-TextCell 'H12' '1'
Set-TextCell 'I12' '12.95'
$ws.Range('J12').Value = 12.95
Set-TextCell 'K12' '850 %'
$ws.Range('L12').Value = 123.02
$ws.Range('M12').Value = 12.3
$ws.Range('N12').Value = 135.32

# Row 13
$ws.Range('A13').Value = 9
$ws.Range('B13').Value = 'wurth'
$ws.Range('C13').Value = '25-06-06'
$ws.Range('D13').Value = 'Carburettor and throttle valve cleaner THRTLVLVECLNR-400G'
Set-TextCell 'E13' '1893100'
Set-TextCell 'F13' '11.21'
Set-TextCell 'G13' '1'
Set-TextCell 'H13' '1'
Set-TextCell 'I13' '11.21'
$ws.Range('J13').Value = 11.21
Set-TextCell 'K13' '850 %'
$ws.Range('L13').Value = 106.5
$ws.Range('M13').Value = 10.65
$ws.Range('N13').Value = 117.15

# Row 14
$ws.Range('A14').Value = 10
$ws.Range('B14').Value = 'wurth'
$ws.Range('C14').Value = '25-06-06'
$ws.Range('D14').Value = 'Deodoriser DEODOR-500ML'
Set-TextCell 'E14' '089313920'
Set-TextCell 'F14' '34.16'
Set-TextCell 'G14' '2'
Set-TextCell 'H14' '2'
Set-TextCell 'I14' '68.32'
$ws.Range('J14').Value = 34.16
Set-TextCell 'K14' '850 %'
$ws.Range('L14').Value = 324.52
$ws.Range('M14').Value = 32.45
$ws.Range('N14').Value = 356.97

# Row 15
$ws.Range('A15').Value = 11
$ws.Range('B15').Value = 'wurth'
$ws.Range('C15').Value = '25-06-06'
$ws.Range('D15').Value = 'Handling FREIGHT LUMP SUM'
Set-TextCell 'E15' '0996000110'
Set-TextCell 'F15' '14.95'
Set-TextCell 'G15' '1'
Set-TextCell 'H15' '1'
Set-TextCell 'I15' '14.95'
$ws.Range('J15').Value = 14.95
Set-TextCell 'K15' '850 %'
$ws.Range('L15').Value = 142.02
$ws.Range('M15').Value = 14.2
$ws.Range('N15').Value = 156.22

